$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'68.349.28"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.359.19"
$ws.Range("E3").Value = "  +0.89%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'585.10"
$ws.Range("E5").Value = "  +0.60%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'177.55"
$ws.Range("E6").Value = "  +0.84%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.03%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.185"
$ws.Range("E9").Value = "  +3.77%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "'0.582"
$ws.Range("E10").Value = "  +0.97%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'48.04"
$ws.Range("E11").Value = "  +5.91%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "'0.0000274"
$ws.Range("E12").Value = "  +2.06%  "

# Row 13 - BitcoinCash
$ws.Range("D13").Value = "'687.68"
$ws.Range("E13").Value = "  +1.97%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'3.909.67"
$ws.Range("E14").Value = "  +0.89%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'8.45"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "'68.390.11"
$ws.Range("E16").Value = "  +1.09%  "

# Row 17 - was TRON, now WrappedEther
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.385.72"
$ws.Range("E17").Value = "  +1.56%  "

# Row 18 - was WrappedEther, now TRON
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.120"
$ws.Range("E18").Value = "  +1.38%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'17.51"
$ws.Range("E19").Value = "  +0.90%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'11.23"
$ws.Range("E20").Value = "  +2.53%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "'0.897"
$ws.Range("E21").Value = "  +0.97%  "

# Row 22 - Toncoin
$ws.Range("E22").Value = "  +0.27%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'16.98"
$ws.Range("E23").Value = "  -0.75%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'100.31"
$ws.Range("E24").Value = "  +1.11%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.89%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "'2.70"
$ws.Range("E26").Value = "  +1.36%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'9.52"
$ws.Range("E27").Value = "  +2.90%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'33.06"
$ws.Range("E28").Value = "  -1.90%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'8.54"
$ws.Range("E29").Value = "  +1.38%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -4.36%  "

# Row 31 - Cosmos
$ws.Range("E31").Value = "  +1.29%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'552.30"

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.105"
$ws.Range("E33").Value = "  +0.69%  "

# Row 34 - OKB
$ws.Range("D34").Value = "'58.14"
$ws.Range("E34").Value = "  +2.65%  "

# Row 35 - Maker
$ws.Range("D35").Value = "'3.721.50"
$ws.Range("E35").Value = "  +0.73%  "

# Row 36 - Dai
$ws.Range("E36").Value = "  +0.12%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  +0.49%  "

# Row 38 - Kaspa
$ws.Range("D38").Value = "'0.138"
$ws.Range("E38").Value = "  +4.72%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = "'34.70"
$ws.Range("E39").Value = "  +1.22%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +2.33%  "

# Row 41 - Fetch.AI
$ws.Range("D41").Value = "'2.62"
$ws.Range("E41").Value = "  +0.06%  "

# Row 42 - PEPE
$ws.Range("D42").Value = "'0.0" + [char]0x2083 + "0675"
$ws.Range("E42").Value = "  +1.22%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  +0.88%  "

# Row 44 - ApeXProtocol
$ws.Range("D44").Value = "'3.26"
$ws.Range("E44").Value = "  -1.16%  "

# Row 45 - VeChain
$ws.Range("D45").Value = "'0.0413"
$ws.Range("E45").Value = "  +1.91%  "

# Row 46 - ThetaToken
$ws.Range("E46").Value = "  +1.64%  "

# Row 47 - Stellar
$ws.Range("D47").Value = "'0.128"
$ws.Range("E47").Value = "  +0.44%  "

# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  -0.01%  "

# Row 49 - Mantle
$ws.Range("D49").Value = "'1.35"
$ws.Range("E49").Value = "  -0.22%  "

# Row 50 - Monero
$ws.Range("D50").Value = "'131.82"
$ws.Range("E50").Value = "  +2.01%  "

# Row 51 - CoreDAO
$ws.Range("D51").Value = "'2.58"
$ws.Range("E51").Value = "  -1.05%  "
